# Auto-generated script applying scheduled market-data refresh to Durandal_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 488.7143
$ws.Range("I18").Value = 294.2
$ws.Range("J18").Value = 975
$ws.Range("K18").Value = 294.2
$ws.Range("L18").Value = 975
$ws.Range("M18").Value = -10.19999999999999
$ws.Range("N18").Value = -1543

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H100").Value = 10754423
$ws.Range("I100").Value = 14493840
$ws.Range("J100").Value = 3599
$ws.Range("K100").Value = 14493840
$ws.Range("L100").Value = 3599
$ws.Range("M100").Value = -14493299
$ws.Range("N100").Value = -4681

$ws.Range("H113").Value = 2177.5557
$ws.Range("I113").Value = 2164
$ws.Range("J113").Value = 2216.2856
$ws.Range("K113").Value = 2164
$ws.Range("L113").Value = 2216.2856
$ws.Range("M113").Value = 1090
$ws.Range("N113").Value = -8724.285599999999

$ws.Range("H138").Value = 4419.492
$ws.Range("I138").Value = 1390.9445
$ws.Range("J138").Value = 8457.556
$ws.Range("K138").Value = 4172.833500000001
$ws.Range("L138").Value = 25372.668
$ws.Range("M138").Value = 967.1664999999994
$ws.Range("N138").Value = -35652.66800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H32").Value = 355877.75
$ws.Range("I32").Value = 2758.4866
$ws.Range("K32").Value = 2758.4866
$ws.Range("M32").Value = -2471.4866

$ws.Range("H74").Value = 808.93335
$ws.Range("I74").Value = 733.38464
$ws.Range("K74").Value = 733.38464
$ws.Range("M74").Value = 140.61536

$ws.Range("H77").Value = 808.93335
$ws.Range("I77").Value = 733.38464
$ws.Range("K77").Value = 3666.9232
$ws.Range("M77").Value = 701.0767999999998

$ws.Range("H119").Value = 31000
$ws.Range("J119").Value = 31000
$ws.Range("L119").Value = 31000
$ws.Range("N119").Value = -40676

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 471
$ws.Range("I22").Value = 321.875
$ws.Range("J22").Value = 769.25
$ws.Range("K22").Value = 321.875
$ws.Range("L22").Value = 769.25
$ws.Range("M22").Value = 28.125
$ws.Range("N22").Value = -1469.25

$ws.Range("H132").Value = 38031.25
$ws.Range("I132").Value = 1653.2941
$ws.Range("J132").Value = 94251.73
$ws.Range("K132").Value = 4959.8823
$ws.Range("L132").Value = 282755.19
$ws.Range("M132").Value = -2429.8823
$ws.Range("N132").Value = -287815.19

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 767.7646999999999
$ws.Range("I5").Value = 924.36365
$ws.Range("J5").Value = 480.66666
$ws.Range("K5").Value = 2773.09095
$ws.Range("L5").Value = 1441.99998
$ws.Range("M5").Value = -2661.09095
$ws.Range("N5").Value = -1665.99998

$ws.Range("H113").Value = 967.98303
$ws.Range("I113").Value = 761.4211
$ws.Range("K113").Value = 2284.2633
$ws.Range("M113").Value = -114.2633000000001

$ws.Range("H122").Value = 618
$ws.Range("I122").Value = 252
$ws.Range("J122").Value = 801
$ws.Range("K122").Value = 2268
$ws.Range("L122").Value = 7209
$ws.Range("M122").Value = 182
$ws.Range("N122").Value = -12109

$ws.Range("H135").Value = 767.7646999999999
$ws.Range("I135").Value = 924.36365
$ws.Range("J135").Value = 480.66666
$ws.Range("K135").Value = 8319.272849999999
$ws.Range("L135").Value = 4325.99994
$ws.Range("M135").Value = -5784.272849999999
$ws.Range("N135").Value = -9395.99994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 10000000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H30").Value = 10000000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H102").Value = 1630.8235
$ws.Range("I102").Value = 1647.3334
$ws.Range("K102").Value = 1647.3334
$ws.Range("M102").Value = -25.33339999999998

$ws.Range("H122").Value = 2027.4546
$ws.Range("I122").Value = 1949.25
$ws.Range("J122").Value = 2236
$ws.Range("K122").Value = 5847.75
$ws.Range("L122").Value = 6708
$ws.Range("M122").Value = -3397.75
$ws.Range("N122").Value = -11608

$ws.Range("H126").Value = 83334340
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 166666670
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 500000010
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -500004950

$ws.Range("H132").Value = 38768.965
$ws.Range("I132").Value = 7877.3125
$ws.Range("J132").Value = 74073.71000000001
$ws.Range("K132").Value = 23631.9375
$ws.Range("L132").Value = 222221.13
$ws.Range("M132").Value = -21101.9375
$ws.Range("N132").Value = -227281.13

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2550.8
$ws.Range("I7").Value = 1956.619
$ws.Range("J7").Value = 3937.2222
$ws.Range("K7").Value = 1956.619
$ws.Range("L7").Value = 3937.2222
$ws.Range("M7").Value = -1844.619
$ws.Range("N7").Value = -4161.2222

$ws.Range("H40").Value = 1701.7693
$ws.Range("I40").Value = 1597.3
$ws.Range("J40").Value = 2050
$ws.Range("K40").Value = 1597.3
$ws.Range("L40").Value = 2050
$ws.Range("M40").Value = -1461.3
$ws.Range("N40").Value = -2322

$ws.Range("H45").Value = 5041
$ws.Range("I45").Value = 5041
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 5041
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4634
$ws.Range("N45").ClearContents()

$ws.Range("H122").Value = 1947.2858
$ws.Range("I122").Value = 1879.8182
$ws.Range("J122").Value = 2061.4614
$ws.Range("K122").Value = 5639.4546
$ws.Range("L122").Value = 6184.3842
$ws.Range("M122").Value = -3189.4546
$ws.Range("N122").Value = -11084.3842

$ws.Range("H126").Value = 2550.8
$ws.Range("I126").Value = 1956.619
$ws.Range("J126").Value = 3937.2222
$ws.Range("K126").Value = 5869.857
$ws.Range("L126").Value = 11811.6666
$ws.Range("M126").Value = -3399.857
$ws.Range("N126").Value = -16751.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12094.25
$ws.Range("I41").Value = 6000
$ws.Range("K41").Value = 6000
$ws.Range("M41").Value = -5610

$ws.Range("H122").Value = 2144.8823
$ws.Range("I122").Value = 2175.6155
$ws.Range("J122").Value = 2045
$ws.Range("K122").Value = 6526.8465
$ws.Range("L122").Value = 6135
$ws.Range("M122").Value = -4076.8465
$ws.Range("N122").Value = -11035

$ws.Range("H126").Value = 41668230
$ws.Range("I126").Value = 76924480
$ws.Range("J126").Value = 1757.7273
$ws.Range("K126").Value = 230773440
$ws.Range("L126").Value = 5273.1819
$ws.Range("M126").Value = -230770970
$ws.Range("N126").Value = -10213.1819

$ws.Range("H138").Value = 79466.664
$ws.Range("J138").Value = 79466.664
$ws.Range("L138").Value = 79466.664
$ws.Range("N138").Value = -89746.664
